$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# "Información SKU" sheet: rename headers.
#   A1: Clave            -> ID Producto
#   F1: Volumen x Pieza  -> Volumen x Unidad
#   G1: Peso x Pieza     -> Peso x Unidad
#   H1: Piezas x Caja    -> Unidades x Caja
#   I1 <-> J1 swap: Peso x Caja / Volumen x Caja
# ----------------------------------------------------------------------
$wsSku = $wb.Worksheets.Item("Información SKU")
$wsSku.Range("A1").Value = "ID Producto"
$wsSku.Range("F1").Value = "Volumen x Unidad"
$wsSku.Range("G1").Value = "Peso x Unidad"
$wsSku.Range("H1").Value = "Unidades x Caja"
$wsSku.Range("I1").Value = "Volumen x Caja"
$wsSku.Range("J1").Value = "Peso x Caja"

# ----------------------------------------------------------------------
# Restore view state (selection/active cell) per sheet, and re-activate
# sheets in the order that leaves "Base de Devoluciones" as the final
# active tab (matching the saved workbook view).
# ----------------------------------------------------------------------

# Foto de Inventarios
$wsFoto = $wb.Worksheets.Item("Foto de Inventarios")
[void]$wsFoto.Activate()
[void]$wsFoto.Range("D38").Select()

# Información SKU
[void]$wsSku.Activate()
[void]$wsSku.Range("L21").Select()

# Base de Recibo
$wsRecibo = $wb.Worksheets.Item("Base de Recibo")
[void]$wsRecibo.Activate()
[void]$wsRecibo.Range("B38").Select()

# Base de Embarque
$wsEmbarque = $wb.Worksheets.Item("Base de Embarque")
[void]$wsEmbarque.Activate()
[void]$wsEmbarque.Range("C11").Select()

# Base de Devoluciones (ends up the active sheet/tab)
$wsDevol = $wb.Worksheets.Item("Base de Devoluciones")
[void]$wsDevol.Activate()
[void]$wsDevol.Range("O13").Select()
